$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: insert a new "Meta description" paragraph right after the title
# (Heading1) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaRange.InsertAfter(": Read our review of Double Diamond, a classic online slot game by IGT. Play for free and enjoy the traditional slot machine experience.")
$metaRange.Collapse(1)
$metaRange.InsertAfter("Meta description")
$metaRange.Bold = 1

# ---------------------------------------------------------------------------
# Part 2: remove the duplicated bold title paragraph near the end of the
# document, and rewrite the following italic paragraph's text with the
# DALLE image prompt (keeping its italic formatting).
# ---------------------------------------------------------------------------
$dupTitlePara = $d.Paragraphs($d.Paragraphs.Count - 1)
$dupTitlePara.Range.Delete()

$promoPara = $d.Paragraphs($d.Paragraphs.Count)

# Find.Execute's replacement text always gets "smart quotes" applied to any
# literal " character, even with AutoFormat options turned off, so we swap
# in a placeholder character for the quotes during the replace, then patch
# the placeholder back to a plain straight quote with a direct Range.Text
# assignment (which does not go through the autocorrect pipeline).
$placeholder = [char]1
$newText = "Prompt: DALLE, please create a cartoon-style feature image for the game " + $placeholder + "Double Diamond" + $placeholder + ". The image should feature a happy Maya warrior with glasses. Make sure to incorporate the Double Diamond logo with the bright gems on either side and include the traditional symbols from old-fashioned mechanical slot machines, such as the numbers 7 and cherries. The image should be eye-catching and engaging to attract potential players."

$promoPara.Range.Find.Execute("Read our review of Double Diamond, a classic online slot game by IGT. Play for free and enjoy the traditional slot machine experience.", $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

$found = $true
while ($found) {
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute($placeholder, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $searchRange.Text = '"'
    }
}
